# Apply the "Fixed update to excel issue" change:
#  1. Rename header B1 on "Weekly Quantity" from "Requested quantity" to "Weekly_PO_Qty"
#  2. Rename header B1 on "Monthly Trend" from "Requested quantity" to "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1 & 2: rename the "Requested quantity" headers ---
$wsWeekly.Cells.Item(1, 2).Value = "Weekly_PO_Qty"
$wsMonthly.Cells.Item(1, 2).Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet at the end ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Copy cell formatting from the existing sheets so styles (bold header,
# date number format) match what's used elsewhere in the workbook.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Data rows
$data = @(
    @(45550.99999999999, 472, 107.630040645277, 800.70843206364),
    @(45578.99999999999, 423, 87.62919905006029, 743.6688489797469),
    @(45585.99999999999, 410, 68.43309191727242, 748.2781732808929),
    @(45592.99999999999, 398, 64.682131647162, 722.3538507873468),
    @(45599.99999999999, 385, 44.85233862956493, 738.1651156896676),
    @(45606.99999999999, 373, 54.23605409404389, 709.1469036455978),
    @(45613.99999999999, 360, 26.16538624086269, 712.6075622530807),
    @(45620.99999999999, 348, 28.3507255107553, 696.740691143207),
    @(45627.99999999999, 336, 5.467552587822535, 646.2551161997754),
    @(45634.99999999999, 323, -24.36403132901136, 652.7789601961937),
    @(45641.99999999999, 311, -18.68027815408103, 630.9080703663409),
    @(45648.99999999999, 298, -53.72815986638719, 656.5344972013704),
    @(45655.99999999999, 286, -59.40935393654546, 626.6903174202349)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
